# Add the new "2022-Q3" quarterly data sheet and update the "总计" (total)
# summary sheet with the new quarter's aggregate row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q3" worksheet (content first, then reposition it so
#    that the freshly created COM reference stays valid while we populate it).
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q3"

# Match the page margins used by the rest of the quarterly sheets
# (0.75in/1in/0.5in - PageSetup margins are expressed in points).
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Header row (row 1), columns B:H - matches the other quarterly sheets.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $c = $newSheet.Cells.Item(1, $i + 2)
    $c.Style = "Heading 1"
    $c.Value = $headers[$i]
}

# Data rows. Column B (fund code) and columns D:G hold text-like numeric
# strings in the source data (e.g. "011686", "0.57"), so force the number
# format to Text before assigning them, otherwise Excel will silently
# reinterpret them as numbers (and "011686" would lose its leading zero).
$newSheet.Range("B2:B3").NumberFormat = "@"
$newSheet.Range("D2:G3").NumberFormat = "@"

# Row 2: 011686 - 创金合信先进装备股票C
$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(2, 2).Value = "011686"
$newSheet.Cells.Item(2, 3).Value = "创金合信先进装备股票C"
$newSheet.Cells.Item(2, 4).Value = "0.57"
$newSheet.Cells.Item(2, 5).Value = "80.17"
$newSheet.Cells.Item(2, 6).Value = "9.92"
$newSheet.Cells.Item(2, 7).Value = "0.0565"
$newSheet.Cells.Item(2, 8).Value = 1

# Row 3: 011685 - 创金合信先进装备股票A
$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Cells.Item(3, 2).Value = "011685"
$newSheet.Cells.Item(3, 3).Value = "创金合信先进装备股票A"
$newSheet.Cells.Item(3, 4).Value = "0.25"
$newSheet.Cells.Item(3, 5).Value = "80.17"
$newSheet.Cells.Item(3, 6).Value = "9.92"
$newSheet.Cells.Item(3, 7).Value = "0.0248"
$newSheet.Cells.Item(3, 8).Value = 1

# Give column A the same style as the other quarterly sheets' index column.
$newSheet.Range("A2:A3").Style = "Heading 1"

# Reposition: right after "总计" and before "2022-Q2" (i.e. position 2).
$newSheet.Move($wb.Worksheets.Item(2))

# ---------------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new row for 2022-Q3 above the existing
#    2022-Q2 row, then renumber the leading index column (A).
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 2
$summary.Cells.Item(2, 4).Value = 0.08
$summary.Range("B2:D2").Style = $summary.Range("B3:D3").Style

# Recompute the 0-based index column now that a row was inserted.
$lastRow = $summary.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $summary.Cells.Item($r, 1)
    $cell.Value = $r - 2
    $cell.Style = "Heading 1"
}
